$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.63
$ws.Range("I3").Value = 2.63
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 3.4
$ws.Range("Z3").Value = 26
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AG3").Value = 7.5
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 10
$ws.Range("AJ3").Value = 26
$ws.Range("AV3").Value = 51
$ws.Range("AW3").Value = 4.5
$ws.Range("AX3").Value = 15
$ws.Range("AY3").Value = 26
$ws.Range("BB3").Value = 201
